# SSDM-13256: remove the "Version" column from the vocabulary terms table.
# The table header row (row 4) used to be: Version | Code | Label | Description
# It becomes:                                Code | Label | Description
# and the two data rows (5 and 6), which both had "1" in the Version column,
# lose that column too, shifting Code/Label/Description one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the content of columns B:D (Code, Label, Description) one column to
# the left, into A:C, for the header + 2 data rows of the terms table.
for ($r = 4; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r, 4).Value2
}

# Remove the now-empty trailing column D entirely so the sheet dimension
# shrinks back down to A:C.
$ws.Columns.Item(4).Delete() | Out-Null

# Rows 4-6 get a slightly smaller row height (15 instead of 15.75) after the
# edit.
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15

# Update the current selection to match the post-edit state.
$ws.Range("A4").Select() | Out-Null
